$d = $word.ActiveDocument

function Split-Run {
    param(
        [string]$SearchText,
        [string]$InnerXml
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $SearchText"
    }
    $s = $rng.Start
    $rng.Delete()
    $insRng = $d.Range($s, $s)
    $xml = @"
<?xml version="1.0" encoding="UTF-16" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            $InnerXml
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
    $insRng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) "Description of some representative Real World Cases  "
#    split into 4 runs, with a gramStart/gramEnd proofErr pair around
#    "Real World"
# ---------------------------------------------------------------------------
Split-Run "Description of some representative Real World Cases  " @'
<w:r><w:t xml:space="preserve">Description of some </w:t></w:r>
<w:r><w:t xml:space="preserve">representative </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>Real World</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> Cases  </w:t></w:r>
'@

# ---------------------------------------------------------------------------
# 2) "Conformance resources for the representation of some study and subject
#     level metadata" -> "Conformance resources and associated examples for
#     the representation of some study and subject level metadata"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Conformance resources for the representation of some study and subject level metadata", $true, $false, $false, $false, $false, $true, 1, $false, "Conformance resources and associated examples for the representation of some study and subject level metadata", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove paragraph "Examples associated to the representative Real World
#    Cases" entirely
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Examples associated to the representative Real World Cases", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph to delete not found" }
$p = $rng.Paragraphs(1)
$p.Range.Delete()

# ---------------------------------------------------------------------------
# 4) "FHIR representation of domain / use case specific data (e.g. genomics;
#     cancer;..)." split into 3 runs with gramStart/gramEnd around
#     "cancer;.."
# ---------------------------------------------------------------------------
Split-Run "FHIR representation of domain / use case specific data (e.g. genomics; cancer;..)." @'
<w:r><w:t xml:space="preserve">FHIR representation of domain / use case specific data (e.g. genomics; </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>cancer;..</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>).</w:t></w:r>
'@

# ---------------------------------------------------------------------------
# 5) "Security protocols (e.g. Data is accessible through an access protocol
#     that supports authentication and authorization)" -> "Privacy and
#     Security management"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Security protocols (e.g. Data is accessible through an access protocol that supports authentication and authorization)", $true, $false, $false, $false, $false, $true, 1, $false, "Privacy and Security management", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) "Consent/Security Labels management" -> "End point" (gramStart/gramEnd)
#    + " discovery"
# ---------------------------------------------------------------------------
Split-Run "Consent/Security Labels management" @'
<w:proofErr w:type="gramStart"/>
<w:r><w:t>End point</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> discovery</w:t></w:r>
'@

# ---------------------------------------------------------------------------
# 7) "FHIR end point discovery (findability)" -> "Provenance management"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("FHIR end point discovery (findability)", $true, $false, $false, $false, $false, $true, 1, $false, "Provenance management", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) Remove paragraph "Provenance resource profiling (Metadata includes
#    provenance information according to community-specific standards)"
#    entirely
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Provenance resource profiling (Metadata includes provenance information according to community-specific standards)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph to delete not found" }
$p = $rng.Paragraphs(1)
$p.Range.Delete()

# ---------------------------------------------------------------------------
# 9) "To be further discussed " -> "Open/Known Issues" with new run
#    formatting
# ---------------------------------------------------------------------------
Split-Run "To be further discussed " @'
<w:r>
  <w:rPr>
    <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
    <w:color w:val="172B4D"/>
    <w:sz w:val="21"/>
    <w:szCs w:val="21"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
  </w:rPr>
  <w:t>Open/Known Issues</w:t>
</w:r>
'@

# ---------------------------------------------------------------------------
# 10) "Persistency of identifiers,  e.g. Digital Object Identifier (DOI) or "
#     split/rearranged; "FEvIR" spellStart/spellEnd + " Object Identifier
#     (FOI)" preserved unchanged after it.
# ---------------------------------------------------------------------------
Split-Run "Persistency of identifiers,  e.g. Digital Object Identifier (DOI) or " @'
<w:r><w:t xml:space="preserve">Persistency of </w:t></w:r>
<w:r><w:t>identifiers</w:t></w:r>
<w:r><w:t xml:space="preserve">, </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>e.g.</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> Digital Object Identifier (DOI) or </w:t></w:r>
'@

# ---------------------------------------------------------------------------
# 11) Remove paragraph '"Metadata is guaranteed to remain available after
#     data is no longer available"' entirely
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute([char]0x22 + "Metadata is guaranteed to remain available after data is no longer available" + [char]0x22, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph to delete not found" }
$p = $rng.Paragraphs(1)
$p.Range.Delete()

# ---------------------------------------------------------------------------
# 12) '"FAIR-compliant vocabularies"' -> 'support for FAIR-compliant
#      vocabularies'
# ---------------------------------------------------------------------------
$d.Content.Find.Execute([char]0x22 + "FAIR-compliant vocabularies" + [char]0x22, $true, $false, $false, $false, $false, $true, 1, $false, "support for FAIR-compliant vocabularies", 2) | Out-Null

Write-Output "All edits applied"
